# Updates "想去人数" (F column) counters on sheet 1 (展览) and sheet 4 (全部类型),
# and inserts a newly-scraped exhibition row (上饶·囧喵喵国风动漫展) into both
# of those sheets, shifting the trailing rows down by one.

$wb = $excel.ActiveWorkbook

function Update-FCounts($ws, $updates) {
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

function Insert-NewRow($ws, $rowNum) {
    # Push rowNum..end down by one row.
    $ws.Rows.Item($rowNum).Insert()

    # Column A uses the bordered/bold/centered header-ish style used by every
    # other row in the "序号" column; copy it down from the row above so the
    # new row matches instead of keeping the blank default style.
    $ws.Range("A" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Range("A" + $rowNum).Value = $rowNum - 1

    # "B" holds a plain date-looking string (not a real Excel date) in this
    # workbook, so force text so COM doesn't silently coerce it to a serial
    # date number.
    $ws.Range("B" + $rowNum).NumberFormat = "@"
    $ws.Range("B" + $rowNum).Value = "2024-08-05"

    $ws.Range("C" + $rowNum).Value = "上饶·囧喵喵国风动漫展"
    $ws.Range("D" + $rowNum).Value = "凤凰东大道与吴楚大道交叉路口北侧 饶派数字文创"
    $ws.Range("E" + $rowNum).Value = "2024.08.05 09:30-08.06 17:30"
    $ws.Range("F" + $rowNum).Value = 8
    $ws.Range("G" + $rowNum).Value = 65
    $ws.Range("H" + $rowNum).Value = "https://show.bilibili.com/platform/detail.html?id=88050"
    $ws.Range("I" + $rowNum).Value = "//i1.hdslb.com/bfs/openplatform/202406/pTwNU5d41719229785696.jpeg"
}

# ---------------------------------------------------------------------------
# Sheet 1 — 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$f1 = @{
    2  = 348
    3  = 279
    4  = 1255
    7  = 60
    9  = 145
    10 = 3476
    11 = 131
    16 = 599
    17 = 94
    18 = 745
    19 = 208
    20 = 122
    22 = 58
    24 = 2628
    25 = 5134
    29 = 3067
    30 = 284
    31 = 2247
    32 = 569
    34 = 83
    35 = 117
    36 = 177
    37 = 313
}
Update-FCounts $ws1 $f1
Insert-NewRow $ws1 38

# ---------------------------------------------------------------------------
# Sheet 4 — 全部类型
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$f4 = @{
    2  = 348
    3  = 279
    4  = 1255
    7  = 60
    9  = 145
    10 = 3476
    11 = 131
    17 = 599
    18 = 94
    19 = 745
    20 = 208
    21 = 122
    23 = 58
    25 = 2628
    26 = 5134
    30 = 3067
    31 = 284
    32 = 2247
    33 = 569
    35 = 83
    36 = 117
    37 = 177
    38 = 313
}
Update-FCounts $ws4 $f4
Insert-NewRow $ws4 39
